$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2025-03-15 Saturday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-03-16 Sunday", 2)

# Update the division problems in the table, cell by cell (row, column) so
# that duplicate expressions (e.g. "32÷5=") are each replaced with the
# correct distinct new value.
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="26÷2="; New="34÷4="},
    @{Row=1;  Col=2; Old="54÷6="; New="84÷6="},
    @{Row=1;  Col=3; Old="70÷9="; New="73÷5="},
    @{Row=1;  Col=4; Old="37÷2="; New="28÷9="},
    @{Row=1;  Col=5; Old="87÷7="; New="17÷9="},

    @{Row=5;  Col=1; Old="93÷7="; New="57÷2="},
    @{Row=5;  Col=2; Old="71÷3="; New="45÷2="},
    @{Row=5;  Col=3; Old="28÷6="; New="22÷7="},
    @{Row=5;  Col=4; Old="42÷2="; New="72÷7="},
    @{Row=5;  Col=5; Old="35÷6="; New="12÷4="},

    @{Row=9;  Col=1; Old="61÷3="; New="60÷6="},
    @{Row=9;  Col=2; Old="27÷7="; New="13÷9="},
    @{Row=9;  Col=3; Old="53÷2="; New="43÷8="},
    @{Row=9;  Col=4; Old="98÷3="; New="78÷9="},
    @{Row=9;  Col=5; Old="97÷8="; New="11÷4="},

    @{Row=13; Col=1; Old="27÷6="; New="35÷9="},
    @{Row=13; Col=2; Old="75÷7="; New="57÷9="},
    @{Row=13; Col=3; Old="68÷4="; New="19÷6="},
    @{Row=13; Col=4; Old="32÷5="; New="23÷4="},
    @{Row=13; Col=5; Old="32÷5="; New="58÷7="},

    @{Row=17; Col=1; Old="24÷4="; New="51÷5="},
    @{Row=17; Col=2; Old="90÷9="; New="60÷6="},
    @{Row=17; Col=3; Old="57÷8="; New="65÷3="},
    @{Row=17; Col=4; Old="44÷8="; New="63÷4="},
    @{Row=17; Col=5; Old="62÷9="; New="74÷6="}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    # Use wdReplaceOne (1) -- not wdReplaceAll -- so that a cell-scoped
    # search only ever touches the single match inside that cell, even
    # when the same expression text appears in more than one cell.
    $rng.Find.Execute($item.Old, $false, $false, $false, $false, $false, $true, 1, $false, $item.New, 1)
}
